# "updated HOAc/Ni(110) areas for uptake"
# Rebuild the data table: new header layout (L, H2, Water, HOAc, CO2, CO)
# across six columns (A:F) and fifteen data rows (2:16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "L"
$ws.Range("B1").Value = "H2"
$ws.Range("C1").Value = "Water"
$ws.Range("D1").Value = "HOAc"
$ws.Range("E1").Value = "CO2"
$ws.Range("F1").Value = "CO"

# The two new header cells (E1, F1) need the same bold/bordered/centered
# header formatting ("s=1") the rest of row 1 already carries. Copy the
# existing header format instead of re-building it, so we reuse the
# existing style record rather than minting a new one.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:F1").PasteSpecial(-4122) | Out-Null

# --- Data rows ----------------------------------------------------------
$data = @(
  @(0,    99432864.00967157, 517946.8848767537, 67958.38511778528, 7559296.938035149, 5077998.739689221),
  @(0.12, 123839430.9698689, 392865.6101648219, 1448068.281180765, 16285203.05489251, 13585008.49200085),
  @(1.04, 122436829.073636,  404390.9086209311, 4105228.491097312, 16213893.11565916, 13434788.10757477),
  @(1.5,  98869495.51061505, 654221.0024057372, 87225.60784737561, 11520319.95079422, 7971459.861646389),
  @(1.8,  103458599.3286628, 565537.5498243775, 55055.75984406852, 11301261.36477978, 8063029.98399758),
  @(2.5,  120815990.679148,  591018.4709477338, 140340.5403203352, 14224533.15584257, 10600218.14778144),
  @(3.2,  117327332.7438248, 572378.1142500686, 181001.7859671152, 14824370.35927413, 10838513.31413946),
  @(3.5,  115699576.6901286, 454232.971138049,  414761.8708728735, 17303903.1349009,  12872728.71279531),
  @(3.6,  127940382.8097512, 446974.5344021593, 232070.0411140699, 15542286.01875807, 11977008.39862336),
  @(4.9,  131997818.9242533, 400388.2611462406, 524824.0145622536, 16787424.83439553, 13387017.79211063),
  @(6.3,  128348734.352567,  262250.7352938735, 461113.9672753855, 14597034.15476643, 12280105.25103936),
  @(7.6,  130262217.033169,  312411.6691422397, 1892970.045781083, 16556193.38364927, 13138063.60631272),
  @(8.4,  122040816.6441483, 245106.2026059856, 840643.6479315368, 15570791.38196047, 12903431.32376256),
  @(9.1,  124884999.5788654, 260454.6925466568, 650669.0681010915, 15471984.68909211, 12206096.30309541),
  @(9.6,  130525130.2219138, 300211.2081272746, 1696378.998010783, 16368221.28626128, 13519583.02790464)
)

$row = 2
foreach ($r in $data) {
  $ws.Cells.Item($row, 1).Value = $r[0]
  $ws.Cells.Item($row, 2).Value = $r[1]
  $ws.Cells.Item($row, 3).Value = $r[2]
  $ws.Cells.Item($row, 4).Value = $r[3]
  $ws.Cells.Item($row, 5).Value = $r[4]
  $ws.Cells.Item($row, 6).Value = $r[5]
  $row++
}

# Rows 13-16 are brand new rows; give column A the same header-matching
# "L" value style ("s=1") the rest of column A already carries, by
# copying the format (only) from an existing styled cell in that column.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A13:A16").PasteSpecial(-4122) | Out-Null
